$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.007.39'
$ws.Range('E2').Value = '  -1.76%  '
$ws.Range('D3').Value = '2.012.67'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('E7').Value = '  +0.03%  '
$c = $ws.Range('D8')
$c.Value = "'54.80"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -5.77%  '
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('E11').Value = '  -5.16%  '
$ws.Range('D12').Value = '2.312.47'
$ws.Range('E12').Value = '  -2.84%  '
$c = $ws.Range('D13')
$c.Value = "'14.10"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -4.81%  '
$c = $ws.Range('D14')
$c.Value = "'20.19"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -4.76%  '
$ws.Range('E15').Value = '  -3.25%  '
$c = $ws.Range('D16')
$c.Value = "'5.10"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('D17').Value = '2.015.51'
$ws.Range('E17').Value = '  -2.81%  '
$ws.Range('D18').Value = '36.967.20'
$ws.Range('E18').Value = '  -1.65%  '
$c = $ws.Range('D19')
$c.Value = "'6.17"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('E21').Value = '  -1.48%  '
$c = $ws.Range('D22')
$c.Value = "'222.70"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.04%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('E25').Value = '  -6.82%  '
$c = $ws.Range('D26')
$c.Value = "'166.21"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.97%  '
$c = $ws.Range('D27')
$c.Value = "'9.14"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -7.61%  '
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D29')
$c.Value = "'18.68"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -3.78%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D30')
$c.Value = "'0.123"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.33%  '
$ws.Range('E31').Value = '  -4.08%  '
$c = $ws.Range('D32')
$c.Value = "'4.51"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.40%  '
$c = $ws.Range('D33')
$c.Value = "'0.0611"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.45%  '
$ws.Range('E34').Value = '  -5.08%  '
$c = $ws.Range('D35')
$c.Value = "'2.34"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -7.38%  '
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('E37').Value = '  +0.18%  '
$c = $ws.Range('D38')
$c.Value = "'3.15"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -4.87%  '
$c = $ws.Range('D39')
$c.Value = "'5.29"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').Value = '1.478.32'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('E41').Value = '  -5.07%  '
$c = $ws.Range('D42')
$c.Value = "'94.72"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.66%  '
$c = $ws.Range('D43')
$c.Value = "'0.0912"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -4.67%  '
$c = $ws.Range('D44')
$c.Value = "'16.24"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -3.36%  '
$ws.Range('E45').Value = '  -5.05%  '
$ws.Range('E46').Value = '  -5.75%  '
$ws.Range('E47').Value = '  -3.10%  '
$c = $ws.Range('D48')
$c.Value = "'7.15"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.50%  '
$c = $ws.Range('D49')
$c.Value = "'2.91"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D50').Value = '2.200.86'
$ws.Range('E50').Value = '  -2.80%  '
$c = $ws.Range('D51')
$c.Value = "'44.35"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.25%  '
